$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4113150465972808
$ws.Range("D2").Value = 0.04122572234101796
$ws.Range("E2").Value = 0.1738898040101944
$ws.Range("F2").Value = 1.595332469086827
$ws.Range("G2").Value = 0.002439077633538228
$ws.Range("L2").Value = 0.151871925898881
$ws.Range("O2").Value = 5.764297367248673
$ws.Range("C3").Value = 0.4025426414724507
$ws.Range("D3").Value = 0.0409752187229131
$ws.Range("E3").Value = 0.1697737840432509
$ws.Range("F3").Value = 1.52810265361525
$ws.Range("G3").Value = 0.002445006841072626
$ws.Range("L3").Value = 0.1479145768505319
$ws.Range("O3").Value = 5.534998599706682
$ws.Range("C4").Value = 0.3974091552251764
$ws.Range("D4").Value = 0.04084533152758496
$ws.Range("E4").Value = 0.1673565443794125
$ws.Range("F4").Value = 1.487888719256958
$ws.Range("G4").Value = 0.002448838950579539
$ws.Range("L4").Value = 0.1455834063616166
$ws.Range("O4").Value = 5.398053472804747
$ws.Range("C5").Value = 0.3953804566620249
$ws.Range("D5").Value = 0.04079838412166481
$ws.Range("E5").Value = 0.1663989955514644
$ws.Range("F5").Value = 1.471766925594508
$ws.Range("G5").Value = 0.002450448910912578
$ws.Range("L5").Value = 0.1446580720818957
$ws.Range("O5").Value = 5.343206354397068
$ws.Range("C6").Value = 0.3950474027089399
$ws.Range("D6").Value = 0.04079094892836821
$ws.Range("E6").Value = 0.1662416511366089
$ws.Range("F6").Value = 1.469105898554105
$ws.Range("G6").Value = 0.002450719168961442
$ws.Range("L6").Value = 0.1445059040535739
$ws.Range("O6").Value = 5.334156736101477
$ws.Range("C7").Value = 0.3973815398675811
$ws.Range("D7").Value = 0.04084467419867366
$ws.Range("E7").Value = 0.1673435194300588
$ws.Range("F7").Value = 1.487670222138732
$ws.Range("G7").Value = 0.002448860466971853
$ws.Range("L7").Value = 0.1455708274521186
$ws.Range("O7").Value = 5.397309912077446
$ws.Range("C8").Value = 0.4082376467136442
$ws.Range("D8").Value = 0.04113436324526987
$ws.Range("E8").Value = 0.1724476636962251
$ws.Range("F8").Value = 1.571929136775651
$ws.Range("G8").Value = 0.002441082374753944
$ws.Range("L8").Value = 0.1504868499154455
$ws.Range("O8").Value = 5.684432272078766
$ws.Range("C9").Value = 0.4315490945530769
$ws.Range("D9").Value = 0.04189385585825534
$ws.Range("E9").Value = 0.1833381646392098
$ws.Range("F9").Value = 1.745725099980262
$ws.Range("G9").Value = 0.002427341336586569
$ws.Range("L9").Value = 0.1609184648719122
$ws.Range("O9").Value = 6.278378879144213
$ws.Range("C10").Value = 0.449934901070776
$ws.Range("D10").Value = 0.04257091168715732
$ws.Range("E10").Value = 0.1918894678433531
$ws.Range("F10").Value = 1.87880370553313
$ws.Range("G10").Value = 0.00241815606995382
$ws.Range("L10").Value = 0.1690778548397418
$ws.Range("O10").Value = 6.734186266718609
$ws.Range("C11").Value = 0.4585782732343091
$ws.Range("D11").Value = 0.04290530830628825
$ws.Range("E11").Value = 0.1959019978053647
$ws.Range("F11").Value = 1.940553300038601
$ws.Range("G11").Value = 0.002414172684373321
$ws.Range("L11").Value = 0.1729002115147011
$ws.Range("O11").Value = 6.945899248271871
$ws.Range("C12").Value = 0.4618919340925629
$ws.Range("D12").Value = 0.04303577475405262
$ws.Range("E12").Value = 0.1974392716319073
$ws.Range("F12").Value = 1.964113322163229
$ws.Range("G12").Value = 0.002412692139231848
$ws.Range("L12").Value = 0.1743637640832532
$ws.Range("O12").Value = 7.026707033161301
$ws.Range("C13").Value = 0.4611764656208379
$ws.Range("D13").Value = 0.04300750511698226
$ws.Range("E13").Value = 0.1971073966082457
$ws.Range("F13").Value = 1.959031345373717
$ws.Range("G13").Value = 0.002413009764156981
$ws.Range("L13").Value = 0.1740478421137368
$ws.Range("O13").Value = 7.009275167049168
$ws.Range("C14").Value = 0.4588500740291579
$ws.Range("D14").Value = 0.04291596470515913
$ws.Range("E14").Value = 0.196028112029289
$ws.Range("F14").Value = 1.942488039117393
$ws.Range("G14").Value = 0.002414050321453078
$ws.Range("L14").Value = 0.1730202948354673
$ws.Range("O14").Value = 6.952534540459965
$ws.Range("C15").Value = 0.4574303918062981
$ws.Range("D15").Value = 0.04286039458749968
$ws.Range("E15").Value = 0.195369345228471
$ws.Range("F15").Value = 1.93237788723826
$ws.Range("G15").Value = 0.002414691318154578
$ws.Range("L15").Value = 0.1723929964417721
$ws.Range("O15").Value = 6.917862456505873
$ws.Range("C16").Value = 0.4493756962201019
$ws.Range("D16").Value = 0.04254959283773729
$ws.Range("E16").Value = 0.1916297216247429
$ws.Range("F16").Value = 1.874792821859529
$ws.Range("G16").Value = 0.002418420303860329
$ws.Range("L16").Value = 0.1688302978953971
$ws.Range("O16").Value = 6.720438916590922
$ws.Range("C17").Value = 0.4445062899978893
$ws.Range("D17").Value = 0.04236571691646418
$ws.Range("E17").Value = 0.1893671001002986
$ws.Range("F17").Value = 1.839778358249163
$ws.Range("G17").Value = 0.002420757748856285
$ws.Range("L17").Value = 0.1666731710256926
$ws.Range("O17").Value = 6.600450126405235
$ws.Range("C18").Value = 0.4417318269260875
$ws.Range("D18").Value = 0.04226243875673674
$ws.Range("E18").Value = 0.1880772180152945
$ws.Range("F18").Value = 1.819752836722699
$ws.Range("G18").Value = 0.002422120553378762
$ws.Range("L18").Value = 0.1654428473609784
$ws.Range("O18").Value = 6.531845942677819
$ws.Range("C19").Value = 0.4407969447695734
$ws.Range("D19").Value = 0.0422278956588471
$ws.Range("E19").Value = 0.1876424582821912
$ws.Range("F19").Value = 1.812992018896381
$ws.Range("G19").Value = 0.002422585135804745
$ws.Range("L19").Value = 0.1650280607491368
$ws.Range("O19").Value = 6.508687920264265
$ws.Range("C20").Value = 0.4450219227506977
$ws.Range("D20").Value = 0.04238503360717516
$ws.Range("E20").Value = 0.1896067665180183
$ws.Range("F20").Value = 1.843493903818882
$ws.Range("G20").Value = 0.002420507024275393
$ws.Range("L20").Value = 0.1669017231800183
$ws.Range("O20").Value = 6.613180618654496
$ws.Range("C21").Value = 0.4595322864317666
$ws.Range("D21").Value = 0.04294274785461027
$ws.Range("E21").Value = 0.1963446389377665
$ws.Range("F21").Value = 1.947342390297166
$ws.Range("G21").Value = 0.002413743929296873
$ws.Range("L21").Value = 0.173321671542169
$ws.Range("O21").Value = 6.969183284980318
$ws.Range("C22").Value = 0.4692525077707899
$ws.Range("D22").Value = 0.04332963540511514
$ws.Range("E22").Value = 0.2008521682772013
$ws.Range("F22").Value = 2.016245085770464
$ws.Range("G22").Value = 0.002409486263778384
$ws.Range("L22").Value = 0.1776114721327957
$ws.Range("O22").Value = 7.205566999775215
$ws.Range("C23").Value = 0.4640428353814912
$ws.Range("D23").Value = 0.04312108369589396
$ws.Range("E23").Value = 0.1984368371623404
$ws.Range("F23").Value = 1.979375152919403
$ws.Range("G23").Value = 0.002411743855429273
$ws.Range("L23").Value = 0.1753132561371586
$ws.Range("O23").Value = 7.079061535084975
$ws.Range("C24").Value = 0.4447887272710034
$ws.Range("D24").Value = 0.04237629295154477
$ws.Range("E24").Value = 0.189498379303636
$ws.Range("F24").Value = 1.841813779557327
$ws.Range("G24").Value = 0.002420620317754585
$ws.Range("L24").Value = 0.1667983641250146
$ws.Range("O24").Value = 6.607423983390618
$ws.Range("C25").Value = 0.4250235504886746
$ws.Range("D25").Value = 0.04166765987363874
$ws.Range("E25").Value = 0.180296338334891
$ws.Range("F25").Value = 1.697774433799481
$ws.Range("G25").Value = 0.002430897980138169
$ws.Range("L25").Value = 0.1580104031887259
$ws.Range("O25").Value = 6.114331337941053
